$d = $word.ActiveDocument

# 1. "There are 2 build rake tasks." -> "There are currently 2 build rake tasks."
$d.Content.Find.Execute(
    "There are 2 build rake tasks.", $true, $false, $false, $false, $false,
    $true, 1, $false, "There are currently 2 build rake tasks.", 2) | Out-Null

# 2. "...Note the trailing ? is required." -> "...Note the trailing / is required."
$d.Content.Find.Execute(
    "Note the trailing ? is required.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Note the trailing / is required.", 2) | Out-Null

# 3. "...by the parallel processors. " -> "...by the parallel search record processors. "
$d.Content.Find.Execute(
    "in the build by the parallel processors.", $true, $false, $false, $false, $false,
    $true, 1, $false, "in the build by the parallel search record processors.", 2) | Out-Null

# 4. "The reason for the 2 ways" -> "The reason for the 3 ways"
$d.Content.Find.Execute(
    "The reason for the 2 ways of creating search records", $true, $false, $false, $false, $false,
    $true, 1, $false, "The reason for the 3 ways of creating search records", 2) | Out-Null

# 5. "entries created and then the search records" -> "entries created first and then the search records"
$d.Content.Find.Execute(
    "I prefer to see the entries created and then the search records", $true, $false, $false, $false, $false,
    $true, 1, $false, "I prefer to see the entries created first and then the search records", 2) | Out-Null

# 6. "takes several days and my crash." -> "takes several days and may crash."
$d.Content.Find.Execute(
    "takes several days and my crash.", $true, $false, $false, $false, $false,
    $true, 1, $false, "takes several days and may crash.", 2) | Out-Null

# 7. Replace "One of these is new. " with the long new sentence (bold run -> plain text)
$rng = $d.Content
$rng.Find.Execute("One of these is new. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$rng.Find.Execute("One of these is new")
$rng.Font.Bold = 0
$rng.Text = "It is possible to process just one file, one userid set of files, one alphabet letter set of userid files, many alphabet set of letters for userids. Also it is possible to create a single county. "

# 8. "a-c  or e-h. The first several of these are likely clear. a-c says process all"
#    -> "a or a-c. The first several of these are likely clear. a says process all userid starting with
#        the letter a (upper of lower case), a-c says process all"
$d.Content.Find.Execute(
    "or a-c  or e-h. The first several of these are likely clear. a-c says process all", $true, $false, $false, $false, $false,
    $true, 1, $false, "or a or a-c. The first several of these are likely clear. a says process all userid starting with the letter a (upper of lower case), a-c says process all", 2) | Out-Null

# 9. "that start with a, b and c or e, f, g and h. This allows" -> "that start with a, b and c. This allows"
$d.Content.Find.Execute(
    "that start with a, b and c or e, f, g and h. This allows", $true, $false, $false, $false, $false,
    $true, 1, $false, "that start with a, b and c. This allows", 2) | Out-Null

# 10. "freereg1_ csv_processor" -> "freereg_ csv_processor"
$d.Content.Find.Execute(
    "possible to run the freereg1_ csv_processor", $true, $false, $false, $false, $false,
    $true, 1, $false, "possible to run the freereg_ csv_processor", 2) | Out-Null

# 11. Insert new paragraph "or by direct calls to the freereg_csv_processor(:type,:search_records,:range)"
#     after "rake build:process_freereg1_csv[:type,:search_records,:range]"
$rng2 = $d.Content
$rng2.Find.Execute("rake build:process_freereg1_csv[:type,:search_records,:range]")
$p = $rng2.Paragraphs(1)
$after = $p.Range
$after.Collapse(0)
$after.InsertParagraphAfter()
$after.Collapse(0)
$after.InsertParagraphAfter()
$after.Collapse(0)
$newPar = $after.Paragraphs(1)
$newPar.Range.Text = "or by direct calls to the freereg_csv_processor(:type,:search_records,:range)"
$newPar.Range.Font.Bold = 1
